$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the pay-period label (week 13 -> week 14). Cells H9 (=B9), B27 (=B9),
# H27 (=B27) and B43 (=H27) are formulas that recompute automatically.
$ws.Range("B9").Value = "SEMANA   14  DEL    04      Al   10   DE   ABRIL          2022"

# Update hours/amount figures in the first pay block.
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 1833

# Update the amount in the second pay block.
$ws.Range("K21").Value = 1820

# Update the extras amount in the third pay block.
$ws.Range("E40").Value = 0

# Restore the view: scroll back to the top and select I44.
$ws.Range("I44").Select()
